$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 680, shifting existing rows 680-721 down to 681-722
$ws.Rows.Item(680).Insert()

# Populate the newly inserted row 680 with the new data point.
# A leading apostrophe forces the date-like text to be stored as a literal
# string instead of being auto-converted to a date serial number; clearing
# the format afterwards drops the "quote prefix" style Excel applies so the
# cell ends up as plain unstyled text, matching the rest of the column.
$ws.Range("A680").Value = "'2026/01/19"
$ws.Range("A680").ClearFormats()
$ws.Range("B680").Value = "月"
$ws.Range("C680").Value = 7
$ws.Range("D680").Value = 173
